# TiffinTracker.xlsx edit script
# - Renames the February sheet's title cell from "January 2018" (shared-string
#   reuse bug in the source) to a new shared string "Feburary 2018".
# - Fills in attendance data for Feb 2nd (row 5) and Feb 5th (row 8).
# - Fills in an Advance payment of 320/360 for Jan (Sat/Sun columns, row 37).
# - Adjusts the saved sheet-view window/selection on both tabs.

$wb = $excel.ActiveWorkbook
$wsJan = $wb.Worksheets.Item("January")
$wsFeb = $wb.Worksheets.Item("February")

# --- January sheet: Advance (row 37) gets values for columns K and L ---
$wsJan.Range("K37").Value = 320
$wsJan.Range("L37").Value = 360

# --- February sheet: title text ---
$wsFeb.Range("B2").Value = "Feburary 2018"

# --- February sheet: attendance rows 5 and 8 ---
$row5 = @(1, 1, 1, 0, 0, 0, 1, 1, 0, 1, 1)
$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsFeb.Range($cols[$i] + "5").Value = $row5[$i]
}

$row8 = @(1, 1, 1, 1, 0, 1, 0, 1, 1, 1, 1)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsFeb.Range($cols[$i] + "8").Value = $row8[$i]
}

# --- Re-apply the "Total Tiffins" row as one shared formula so the saved
#     file groups D32:N32 the same way Excel does when the whole row is
#     filled/edited together ---
$wsFeb.Range("D32:N32").Formula = "=SUM(D4:D29)"

# --- Recalculate so dependent formulas refresh their cached values ---
$excel.Calculate()

# --- Restore the saved window/selection state for each sheet ---
$wsJan.Activate()
$wsJan.Application.ActiveWindow.ScrollRow = 22
$wsJan.Range("F16").Select()

$wsFeb.Activate()
$wsFeb.Application.ActiveWindow.ScrollRow = 10
$wsFeb.Range("J21").Select()
